# bookkeeper-report.xlsx fix:
# The training-release labels lost their numeric prefix ("1, 4.1.0" -> ", 4.0.0",
# "2, 4.2.0" -> ", 4.0.0, 4.1.0", etc.) and several derived metric columns
# (Precision/Recall/AUC/Kappa) were recomputed now that the file-path matching
# bug is fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a number need to be pinned to Text format
# first, otherwise Excel will silently convert the numeric-looking string into
# a real numeric cell instead of keeping it as the shared-string text that the
# report expects.
$textCells = @(
    "F2",
    "D5", "E5", "F5", "G5",
    "E6", "F6", "G6",
    "D7", "E7", "F7", "G7",
    "F8",
    "F9",
    "F11",
    "F12",
    "D13", "E13", "F13", "G13"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 (Bookkeeper / Naive Bayes)
$ws.Range("B2").Value = ", 4.0.0"
$ws.Range("F2").Value = "0.29120879120879123"

# Row 3 (Bookkeeper / Random Forest)
$ws.Range("B3").Value = ", 4.0.0"

# Row 4 (Bookkeeper / IBk)
$ws.Range("B4").Value = ", 4.0.0"

# Row 5 (Bookkeeper / Naive Bayes)
$ws.Range("B5").Value = ", 4.0.0, 4.1.0"
$ws.Range("D5").Value = "0.35294117647058826"
$ws.Range("E5").Value = "0.2727272727272727"
$ws.Range("F5").Value = "0.689935064935065"
$ws.Range("G5").Value = "0.20993014914102326"

# Row 6 (Bookkeeper / Random Forest)
$ws.Range("B6").Value = ", 4.0.0, 4.1.0"
$ws.Range("D6").Value = "NaN"
$ws.Range("E6").Value = "0.0"
$ws.Range("F6").Value = "0.6883116883116883"
$ws.Range("G6").Value = "0.0"

# Row 7 (Bookkeeper / IBk)
$ws.Range("B7").Value = ", 4.0.0, 4.1.0"
$ws.Range("D7").Value = "0.25"
$ws.Range("E7").Value = "0.045454545454545456"
$ws.Range("F7").Value = "0.5673701298701299"
$ws.Range("G7").Value = "0.03476907109496633"

# Row 8 (Bookkeeper / Naive Bayes)
$ws.Range("B8").Value = ", 4.0.0, 4.1.0, 4.2.0"
$ws.Range("F8").Value = "0.8030587833219412"

# Row 9 (Bookkeeper / Random Forest)
$ws.Range("B9").Value = ", 4.0.0, 4.1.0, 4.2.0"
$ws.Range("F9").Value = "0.5658749145591251"

# Row 10 (Bookkeeper / IBk)
$ws.Range("B10").Value = ", 4.0.0, 4.1.0, 4.2.0"

# Row 11 (Bookkeeper / Naive Bayes)
$ws.Range("B11").Value = ", 4.0.0, 4.1.0, 4.2.0, 4.2.1"
$ws.Range("F11").Value = "0.8496534965349654"

# Row 12 (Bookkeeper / Random Forest)
$ws.Range("B12").Value = ", 4.0.0, 4.1.0, 4.2.0, 4.2.1"
$ws.Range("F12").Value = "0.7145171451714517"

# Row 13 (Bookkeeper / IBk)
$ws.Range("B13").Value = ", 4.0.0, 4.1.0, 4.2.0, 4.2.1"
$ws.Range("D13").Value = "0.6"
$ws.Range("E13").Value = "0.07317073170731707"
$ws.Range("F13").Value = "0.8068130681306813"
$ws.Range("G13").Value = "0.10486300387318909"
